$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF)
# Copy formatting from the existing header cell (H1) so the new
# header cells share the same cell style (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I and J, rows 2-30
$data = @(
    @(2, 7, 7),
    @(3, 6, 6),
    @(4, 7, 7),
    @(5, 7, 7),
    @(6, 4, 5),
    @(7, 7, 7),
    @(8, 6, 6),
    @(9, 8, 8),
    @(10, 6, 6),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 7, 7),
    @(14, 7, 8),
    @(15, 9, 10),
    @(16, 6, 6),
    @(17, 4, 5),
    @(18, 8, 8),
    @(19, 5, 6),
    @(20, 3, 4),
    @(21, 6, 7),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 4, 6),
    @(25, 8, 8),
    @(26, 5, 6),
    @(27, 7, 8),
    @(28, 3, 3),
    @(29, 7, 7),
    @(30, 8, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
